$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells that receive numeric-looking text
# stay stored as text, matching the inlineStr type used in the workbook.
$priceCells = @("D2","D3","D5","D6","D9","D12","D13","D14","D15","D16","D17","D18","D21","D25","D26","D27","D29","D31","D37","D38","D39","D41","D42","D44","D45","D46","D48","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Cells.Item(2, 4).Value = '68.256.97'

# Row 3
$ws.Cells.Item(3, 4).Value = '2.641.64'
$ws.Cells.Item(3, 5).Value = '  +0.53%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.01%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '598.26'
$ws.Cells.Item(5, 5).Value = '  +0.29%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '154.60'
$ws.Cells.Item(6, 5).Value = '  +0.54%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.01%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.86%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '2.639.59'
$ws.Cells.Item(9, 5).Value = '  +0.51%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +7.62%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -0.76%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '5.26'
$ws.Cells.Item(12, 5).Value = '  +0.80%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '0.354'
$ws.Cells.Item(13, 5).Value = '  +1.91%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '0.0000194'
$ws.Cells.Item(14, 5).Value = '  +2.98%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '27.91'
$ws.Cells.Item(15, 5).Value = '  +1.10%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '3.121.17'
$ws.Cells.Item(16, 5).Value = '  +0.53%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '68.151.56'
$ws.Cells.Item(17, 5).Value = '  +0.63%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '2.643.19'
$ws.Cells.Item(18, 5).Value = '  +0.43%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -0.48%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  -1.33%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '7.45'
$ws.Cells.Item(21, 5).Value = '  +0.22%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +1.20%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -0.96%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '75.26'
$ws.Cells.Item(25, 5).Value = '  +4.53%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '0.999'
$ws.Cells.Item(26, 5).Value = '  -0.16%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '9.75'
$ws.Cells.Item(27, 5).Value = '  -0.90%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +2.09%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '2.776.17'
$ws.Cells.Item(29, 5).Value = '  +0.54%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +0.01%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '562.38'
$ws.Cells.Item(31, 5).Value = '  -2.12%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +1.75%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +0.39%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +0.96%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +2.21%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -0.02%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '1.57'
$ws.Cells.Item(37, 5).Value = '  +2.50%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '161.50'
$ws.Cells.Item(38, 5).Value = '  +1.78%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '19.30'
$ws.Cells.Item(39, 5).Value = '  +0.82%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +1.79%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '1.89'
$ws.Cells.Item(41, 5).Value = '  -0.04%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '5.34'
$ws.Cells.Item(42, 5).Value = '  -0.34%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  +1.23%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'dogwifhat'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(44, 4).Value = '2.64'
$ws.Cells.Item(44, 5).Value = '  +0.25%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(45, 4).Value = '17.75'
$ws.Cells.Item(45, 5).Value = '  +2.10%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '40.63'
$ws.Cells.Item(46, 5).Value = '  +1.45%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -0.09%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '156.07'
$ws.Cells.Item(48, 5).Value = '  +0.45%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +1.75%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +0.12%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(51, 4).Value = '21.78'
$ws.Cells.Item(51, 5).Value = '  -0.85%  '
